# This workbook lists species observations in rows 9-18 (columns A-R hold the
# per-observation data: Id, Taxonsorteringsordning, Rödlistade, TaxonId,
# Artnamn, Vetenskapligt namn, Auktor, Ost, Nord). The edit re-orders which
# observation record sits in which row (the surrounding metadata in columns
# S onward - accuracy, county, municipality, dates, reporters - is identical
# for every one of these rows, so only the per-record columns need to move).
#
# Values are written directly (computed from the known before/after state)
# rather than read-and-rewritten in place, since several rows' contents get
# swapped/rotated simultaneously and an in-place read-after-write could pick
# up an already-overwritten value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 (was row 15: Id 111671406)
$ws.Range("A9").Value = 111671406
$ws.Range("B9").Value = 78578
$ws.Range("D9").Value = "NT"
$ws.Range("E9").Value = 6458
$ws.Range("F9").Value = "Lunglav"
$ws.Range("G9").Value = "Lobaria pulmonaria"
$ws.Range("H9").Value = "(L.) Hoffm."
$ws.Range("Q9").Value = 557823.3030943703
$ws.Range("R9").Value = 7068159.357501161
$ws.Range("L9").ClearContents()

# Row 10 (was row 14: Id 111671364)
$ws.Range("A10").Value = 111671364
$ws.Range("B10").Value = 96368
$ws.Range("D10").Value = "LC"
$ws.Range("E10").Value = 221952
$ws.Range("F10").Value = "Spindelblomster"
$ws.Range("G10").Value = "Neottia cordata"
$ws.Range("H10").Value = "(L.) Rich."
$ws.Range("Q10").Value = 557813.3601359134
$ws.Range("R10").Value = 7068169.364891288

# Row 11 (was row 17: Id 111671384)
$ws.Range("A11").Value = 111671384
$ws.Range("B11").Value = 96348
$ws.Range("D11").Value = "VU"
$ws.Range("E11").Value = 220787
$ws.Range("F11").Value = "Knärot"
$ws.Range("G11").Value = "Goodyera repens"
$ws.Range("H11").Value = "(L.) R. Br."
$ws.Range("Q11").Value = 557798.0632258818
$ws.Range("R11").Value = 7068181.046264404

# Row 12 (was row 10: Id 111671395)
$ws.Range("A12").Value = 111671395
$ws.Range("B12").Value = 96348
$ws.Range("D12").Value = "VU"
$ws.Range("E12").Value = 220787
$ws.Range("F12").Value = "Knärot"
$ws.Range("G12").Value = "Goodyera repens"
$ws.Range("H12").Value = "(L.) R. Br."
$ws.Range("Q12").Value = 557763.2623863788
$ws.Range("R12").Value = 7068264.582601988

# Row 13 (was row 11: Id 111670575)
$ws.Range("A13").Value = 111670575
$ws.Range("B13").Value = 96346
$ws.Range("D13").Value = "NT"
$ws.Range("E13").Value = 620
$ws.Range("F13").Value = "Skogsfru"
$ws.Range("G13").Value = "Epipogium aphyllum"
$ws.Range("H13").Value = "Sw."
$ws.Range("Q13").Value = 558082.6649719321
$ws.Range("R13").Value = 7067974.943554637

# Row 14 (was row 18: Id 111670607)
$ws.Range("A14").Value = 111670607
$ws.Range("B14").Value = 96368
$ws.Range("D14").Value = "LC"
$ws.Range("E14").Value = 221952
$ws.Range("F14").Value = "Spindelblomster"
$ws.Range("G14").Value = "Neottia cordata"
$ws.Range("H14").Value = "(L.) Rich."
$ws.Range("Q14").Value = 558031.5471372061
$ws.Range("R14").Value = 7067907.98648507

# Row 15 (was row 16: Id 111670593)
$ws.Range("A15").Value = 111670593
$ws.Range("B15").Value = 78578
$ws.Range("D15").Value = "NT"
$ws.Range("E15").Value = 6458
$ws.Range("F15").Value = "Lunglav"
$ws.Range("G15").Value = "Lobaria pulmonaria"
$ws.Range("H15").Value = "(L.) Hoffm."
$ws.Range("Q15").Value = 558040.5475534229
$ws.Range("R15").Value = 7067901.063021242

# Row 16 (was row 13: Id 111671345)
$ws.Range("A16").Value = 111671345
$ws.Range("B16").Value = 96348
$ws.Range("D16").Value = "VU"
$ws.Range("E16").Value = 220787
$ws.Range("F16").Value = "Knärot"
$ws.Range("G16").Value = "Goodyera repens"
$ws.Range("H16").Value = "(L.) R. Br."
$ws.Range("Q16").Value = 557812.5300353739
$ws.Range("R16").Value = 7068166.248475613
$ws.Range("L16").Value = ""

# Row 17 (was row 12: Id 111670588)
$ws.Range("A17").Value = 111670588
$ws.Range("B17").Value = 96348
$ws.Range("D17").Value = "VU"
$ws.Range("E17").Value = 220787
$ws.Range("F17").Value = "Knärot"
$ws.Range("G17").Value = "Goodyera repens"
$ws.Range("H17").Value = "(L.) R. Br."
$ws.Range("Q17").Value = 558039.6361001397
$ws.Range("R17").Value = 7067902.375451046

# Row 18 (was row 9: Id 111670599)
$ws.Range("A18").Value = 111670599
$ws.Range("B18").Value = 96348
$ws.Range("D18").Value = "VU"
$ws.Range("E18").Value = 220787
$ws.Range("F18").Value = "Knärot"
$ws.Range("G18").Value = "Goodyera repens"
$ws.Range("H18").Value = "(L.) R. Br."
$ws.Range("Q18").Value = 558031.5226908802
$ws.Range("R18").Value = 7067909.315233406
